$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text in Price/Volume columns stays as text (matches source inlineStr formatting)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.877.67"
$ws.Range("D3").Value = "1.816.10"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "308.35"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("D7").Value = "0.4618"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("D8").Value = "0.3654"
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").Value = "0.07221"
$ws.Range("D10").Value = "0.8577"
$ws.Range("E10").Value = "  -3.22%  "
$ws.Range("D11").Value = "19.69"
$ws.Range("E11").Value = "  -3.83%  "
$ws.Range("D12").Value = "0.07524"
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.734.35"
$ws.Range("E13").Value = "  -7.80%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.319"
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("D15").Value = "91.75"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "6.474"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D18").Value = "0.000008598"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "14.42"
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "26.641.95"
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("D23").Value = "10.50"
$ws.Range("D24").Value = "1.956.82"
$ws.Range("E24").Value = "  -6.63%  "
$ws.Range("D25").Value = "151.73"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  -3.73%  "
$ws.Range("D27").Value = "18.10"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("D28").Value = "2.075"
$ws.Range("E28").Value = "  -3.05%  "
$ws.Range("D29").Value = "5.087"
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("D30").Value = "115.06"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").Value = "0.08863"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "4.408"
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("D34").Value = "1.126"
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("D35").Value = "0.7138"
$ws.Range("E35").Value = "  -5.70%  "
$ws.Range("D36").Value = "1.075"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").Value = "0.05235"
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").Value = "2.409"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").Value = "0.01913"
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("D40").Value = "2.920"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("D41").Value = "7.132"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").Value = "0.5137"
$ws.Range("E42").Value = "  -3.69%  "
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D44").Value = "8.153"
$ws.Range("E44").Value = "  -4.24%  "
$ws.Range("D45").Value = "0.4791"
$ws.Range("E45").Value = "  -2.71%  "
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.10"
$ws.Range("E47").Value = "  -4.55%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "102.87"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").Value = "0.06275"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("D51").Value = "63.78"
$ws.Range("E51").Value = "  -3.08%  "
